$wb = $excel.ActiveWorkbook

# --- Add the new sheet "Repute Traders 103" after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Repute Traders 103"

# --- Header row ---
$ws.Range("A1").Value = "SR. NO"
$ws.Range("B1").Value = "Product Description"
$ws.Range("C1").Value = "Quantity"
$ws.Range("D1").Value = "Pricing Per Unit"
$ws.Range("E1").Value = "Total"

$headerRange = $ws.Range("A1:E1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108
$headerRange.WrapText = $true
$headerRange.Borders.LineStyle = 1

# --- Data row ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Canvas Hoses 101 mm"
$ws.Range("C2").Value = 30.5
$ws.Range("D2").Value = 66
$ws.Range("E2").Formula = "=30.5*66"

$dataRange = $ws.Range("A2:E2")
$dataRange.HorizontalAlignment = -4108
$dataRange.VerticalAlignment = -4108
$dataRange.WrapText = $true
$dataRange.Borders.LineStyle = 1

# --- Total row ---
$ws.Range("A3").Value = "Total"
$ws.Range("E3").Formula = "=SUM(E2)"
$ws.Range("A3:D3").Merge()

# --- GST 12% row ---
$ws.Range("A4").Value = "GST 12%"
$ws.Range("E4").Formula = "=E3*12%"
$ws.Range("A4:D4").Merge()

# --- Grand total row ---
$ws.Range("A5").Value = "Grand Total"
$ws.Range("E5").Formula = "=SUM(E3:E4)"
$ws.Range("A5:D5").Merge()

$totalsRange = $ws.Range("A3:E5")
$totalsRange.Font.Bold = $true
$totalsRange.HorizontalAlignment = -4108
$totalsRange.VerticalAlignment = -4108
$totalsRange.WrapText = $true
$totalsRange.Borders.LineStyle = 1

# Column widths roughly matching the source layout
$ws.Columns.Item(2).ColumnWidth = 26.33
$ws.Columns.Item(3).ColumnWidth = 9.11
$ws.Columns.Item(4).ColumnWidth = 16.66
$ws.Columns.Item(5).ColumnWidth = 6.11

$ws.Activate()
